$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray -1 value in F56 (no longer applicable)
$ws.Range("F56").ClearContents()

# Insert a new row above row 58 for the new component
# "Lithium-Ion-LFP-bicharger % discharge" (Link), pushing existing rows 58+ down by one
$ws.Rows("58:58").Insert()

# Populate the newly inserted row 58 with the new Link component data
$ws.Range("A58").Value = "Link"
$ws.Range("B58").Value = "Lithium-Ion-LFP-bicharger % discharge"
$ws.Range("C58").Value = "battery"
$ws.Range("D58").Value = "lithium_ion_lfp"
$ws.Range("E58").Value = "electricity"
$ws.Range("I58").Value = 0
$ws.Range("J58").Formula = "=B42 & ""/time range/"" & B41"
$ws.Range("K58").Value = "db"
$ws.Range("L58").Formula = "= B42 & ""/"" & B41 & B40"
$ws.Range("N58").Value = "db"
$ws.Range("O58").Value = "db"

# Update the active selection to reflect where the edit was made
[void]$ws.Range("F59").Select()
